# Quarterly indexing esoteric bug-fix operation
#
# Column A (rows 2..last) holds quarter-start dates (serial dates stored as
# the 1st day of a quarter's start month, e.g. 2020-01-01, 2020-04-01, ...).
# The fix re-indexes each one to the 15th of the *following* month
# (e.g. 2020-01-01 -> 2020-02-15), which is how the naive QoQ forecaster is
# meant to timestamp its quarterly observations.
#
# Helper: days-since-1899-12-30 (the Excel date epoch, accounting for the
# fictitious 1900 leap day) computed via an integer civil-calendar algorithm,
# since DateTime subtraction isn't reliably supported in this environment.
function Get-DaysFromCivil($year, $month, $day) {
    $y = $year
    if ($month -le 2) { $y = $y - 1 }
    $m = $month
    $era = [Math]::Floor((if ($y -ge 0) { $y } else { $y - 399 }) / 400)
    $yoe = $y - $era * 400
    $doyShift = [Math]::Floor((153 * (if ($m -gt 2) { $m - 3 } else { $m + 9 }) + 2) / 5)
    $doy = $doyShift + $day - 1
    $doe = $yoe * 365 + [Math]::Floor($yoe / 4) - [Math]::Floor($yoe / 100) + $doy
    $days = $era * 146097 + $doe - 719468
    return $days
}

$excelEpochDays = Get-DaysFromCivil 1899 12 30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $serial = $cell.Value2

    # Convert the existing serial date to a real DateTime so we can inspect
    # its year/month.
    $epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
    $d = $epoch.AddDays($serial)

    # Move forward one month, and pin the day-of-month to 15.
    $newMonth = $d.Month + 1
    $newYear = $d.Year
    if ($newMonth -gt 12) {
        $newMonth = $newMonth - 12
        $newYear = $newYear + 1
    }

    $newDays = Get-DaysFromCivil $newYear $newMonth 15
    $newSerial = $newDays - $excelEpochDays

    $cell.Value2 = $newSerial
}
